# Update "想去人数" (want-to-go count) values on the sheets that hold
# the exhibition data: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 60
    $ws.Range("F5").Value = 113
}
